# feat: add 2022-Q1 data
#
# Workbook currently has two sheets: "2021-Q2" (fund snapshot) and "总计"
# (rollup of all quarters). This adds a new "2022-Q1" fund snapshot sheet
# (cloned from "2021-Q2", with updated numbers) positioned right after
# "2021-Q2", and records the new quarter as the first data row of "总计".

$wb = $excel.ActiveWorkbook

$sheetQ2 = $wb.Worksheets.Item(1)      # "2021-Q2"

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet by duplicating "2021-Q2" (this carries
#    over all formatting/styles) and placing it right after it.
# ---------------------------------------------------------------------
$sheetQ2.Copy($null, $sheetQ2)
$sheetQ1 = $wb.Worksheets.Item("2021-Q2 (2)")
$sheetQ1.Name = "2022-Q1"

# Re-fetch "总计" by name: inserting the new sheet shifts its position,
# and a previously grabbed positional reference would now point at the
# wrong sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row: only the "基金金额" column is renamed to "基金规模".
$sheetQ1.Range("D1").Value = "基金规模"

# Data row 2: fund code/name (B2/C2) stay the same; update the rest.
# D2, E2, F2, G2 are stored as text in the source file, so force text
# with a leading apostrophe (quote-prefix) to avoid Excel auto-converting
# the numeric-looking strings into real numbers.
$sheetQ1.Range("D2").Value = "'0.27"
$sheetQ1.Range("E2").Value = "'89.72"
$sheetQ1.Range("F2").Value = "'1.27"
$sheetQ1.Range("G2").Value = "'0.0034"
$sheetQ1.Range("H2").Value = 3

# ---------------------------------------------------------------------
# 2. Record the new quarter in "总计": insert a new row right under the
#    header and push the existing "2021-Q2" row down.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Restore formatting on the new row's A cell (bold, centered, boxed --
# matching the look of the other index cells in column A) and clear the
# incidental formatting the row-insert applied to B2:D2.
$totalSheet.Range("B2:D2").ClearFormats()
$indexCell = $totalSheet.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0

# The row that used to be index 0 (2021-Q2) is now the second data row.
$totalSheet.Range("A3").Value = 1
